# Delete Application table and move attachments to documents
#
# This reproduces, via Excel COM automation, the following changes to
# "Projekt Dom.xlsx":
#   - On sheet "Formalnosci": the little "POSTS" reference/lookup table that
#     lived next to the "STAGE_I" / "SEND_TYPE" tables (columns J:M, rows
#     12-16) is removed, its "SEND_TYPE" neighbour collapses from F:H down to
#     a single column (E), and the "POSTS" table is recreated lower on the
#     sheet, in column C only, as rows 20-24.
#   - The active sheet/window selection moves: "EtapI " is no longer the
#     selected tab, "Formalnosci" becomes the active tab (with a new
#     selection + scroll position), and the selection on "EtapII" changes.

$wb = $excel.ActiveWorkbook

$wsEtapI       = $wb.Worksheets.Item("EtapI ")
$wsFormalnosci = $wb.Worksheets.Item("Formalnosci")
$wsEtapII      = $wb.Worksheets.Item("EtapII")

# ---------------------------------------------------------------------
# 1. Rework the small lookup tables around C12:M16 on "Formalnosci".
# ---------------------------------------------------------------------

# Remember the "SEND_TYPE" table header text (currently in F12) before we
# start clearing cells, so we can move it over to E12.
$sendTypeHeader = $wsFormalnosci.Range("F12").Value2

# Remove the now-unused "SEND_TYPE" helper cells in columns F:H (rows 12-13)
# and the whole "POSTS" table that used to sit in columns J:M (rows 12-16).
$wsFormalnosci.Range("F12:H13").Clear()
$wsFormalnosci.Range("J12:M16").Clear()

# The "SEND_TYPE" table now only needs a single column: put its header back
# in E12 (re-using the same banner style that was already on that cell).
$wsFormalnosci.Range("E12").Value = $sendTypeHeader

# Re-create the "POSTS" table lower on the sheet (column C, rows 20-24),
# copying the banner/header formatting from the equivalent rows of the
# "STAGE_I" table (C12/C13) so the styling matches.
$wsFormalnosci.Range("C12").Copy()
$wsFormalnosci.Range("C20").PasteSpecial(-4122)
$wsFormalnosci.Range("C20").Value = "POSTS"

$wsFormalnosci.Range("C13").Copy()
$wsFormalnosci.Range("C21").PasteSpecial(-4122)
$wsFormalnosci.Range("C21").Value = "name"

$wsFormalnosci.Range("C22").Value = "starostwo powiatowe"
$wsFormalnosci.Range("C23").Value = "pge"
$wsFormalnosci.Range("C24").Value = "urząd gminy"

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 2. Update the active tab / selections to match the new window state.
# ---------------------------------------------------------------------

# "EtapII" selection moves from E5 to I3 (it stays a background sheet).
$wsEtapII.Range("I3").Select()

# "Formalnosci" becomes the active sheet, scrolled so row 4 is at the top,
# with F12 selected. Activating it also clears "EtapI "'s tabSelected flag.
$wsFormalnosci.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 4
$win.ScrollColumn = 1
$wsFormalnosci.Range("F12").Select()
